{"js": "// Replace the date line and each \"A\u00f7B=C, D\" division answer in the table\n// with the new values from the commit. Every old value is a unique string\n// in the document, so a straightforward exact-text search/replace for each\n// pair is robust and keeps the original run formatting (font/size) intact,\n// since insertText(..., Replace) only swaps the text inside the existing run.\nconst replacements = [\n  [\"2025-09-22 Monday\", \"2025-09-23 Tuesday\"],\n  [\"515\\u00F77=73, 4\", \"592\\u00F79=65, 7\"],\n  [\"424\\u00F72=212, 0\", \"150\\u00F74=37, 2\"],\n  [\"837\\u00F73=279, 0\", \"711\\u00F77=101, 4\"],\n  [\"823\\u00F74=205, 3\", \"824\\u00F77=117, 5\"],\n  [\"530\\u00F77=75, 5\", \"947\\u00F74=236, 3\"],\n  [\"153\\u00F77=21, 6\", \"130\\u00F75=26, 0\"],\n  [\"729\\u00F78=91, 1\", \"763\\u00F72=381, 1\"],\n  [\"856\\u00F78=107, 0\", \"123\\u00F77=17, 4\"],\n  [\"319\\u00F77=45, 4\", \"621\\u00F75=124, 1\"],\n  [\"806\\u00F79=89, 5\", \"716\\u00F75=143, 1\"],\n  [\"962\\u00F78=120, 2\", \"357\\u00F73=119, 0\"],\n  [\"393\\u00F78=49, 1\", \"704\\u00F74=176, 0\"],\n  [\"761\\u00F74=190, 1\", \"340\\u00F72=170, 0\"],\n  [\"307\\u00F78=38, 3\", \"530\\u00F79=58, 8\"],\n  [\"511\\u00F78=63, 7\", \"356\\u00F75=71, 1\"],\n  [\"497\\u00F75=99, 2\", \"159\\u00F73=53, 0\"],\n  [\"845\\u00F73=281, 2\", \"297\\u00F77=42, 3\"],\n  [\"661\\u00F75=132, 1\", \"837\\u00F72=418, 1\"],\n  [\"232\\u00F72=116, 0\", \"846\\u00F73=282, 0\"],\n  [\"280\\u00F72=140, 0\", \"324\\u00F79=36, 0\"],\n  [\"715\\u00F74=178, 3\", \"942\\u00F78=117, 6\"],\n  [\"132\\u00F73=44, 0\", \"923\\u00F79=102, 5\"],\n  [\"233\\u00F76=38, 5\", \"375\\u00F77=53, 4\"],\n  [\"995\\u00F74=248, 3\", \"775\\u00F74=193, 3\"],\n  [\"408\\u00F78=51, 0\", \"227\\u00F75=45, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"A\u00f7B=C, D\" division answer in the table\n# with the new values from the commit. Every old value is a unique string\n# in the document, so a Find/Replace (wdReplaceOne) for each pair is robust\n# and preserves the original run formatting (font/size), because Find just\n# swaps the text inside the matched range.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-09-22 Monday\", \"2025-09-23 Tuesday\"),\n    @(\"515\u00f77=73, 4\", \"592\u00f79=65, 7\"),\n    @(\"424\u00f72=212, 0\", \"150\u00f74=37, 2\"),\n    @(\"837\u00f73=279, 0\", \"711\u00f77=101, 4\"),\n    @(\"823\u00f74=205, 3\", \"824\u00f77=117, 5\"),\n    @(\"530\u00f77=75, 5\", \"947\u00f74=236, 3\"),\n    @(\"153\u00f77=21, 6\", \"130\u00f75=26, 0\"),\n    @(\"729\u00f78=91, 1\", \"763\u00f72=381, 1\"),\n    @(\"856\u00f78=107, 0\", \"123\u00f77=17, 4\"),\n    @(\"319\u00f77=45, 4\", \"621\u00f75=124, 1\"),\n    @(\"806\u00f79=89, 5\", \"716\u00f75=143, 1\"),\n    @(\"962\u00f78=120, 2\", \"357\u00f73=119, 0\"),\n    @(\"393\u00f78=49, 1\", \"704\u00f74=176, 0\"),\n    @(\"761\u00f74=190, 1\", \"340\u00f72=170, 0\"),\n    @(\"307\u00f78=38, 3\", \"530\u00f79=58, 8\"),\n    @(\"511\u00f78=63, 7\", \"356\u00f75=71, 1\"),\n    @(\"497\u00f75=99, 2\", \"159\u00f73=53, 0\"),\n    @(\"845\u00f73=281, 2\", \"297\u00f77=42, 3\"),\n    @(\"661\u00f75=132, 1\", \"837\u00f72=418, 1\"),\n    @(\"232\u00f72=116, 0\", \"846\u00f73=282, 0\"),\n    @(\"280\u00f72=140, 0\", \"324\u00f79=36, 0\"),\n    @(\"715\u00f74=178, 3\", \"942\u00f78=117, 6\"),\n    @(\"132\u00f73=44, 0\", \"923\u00f79=102, 5\"),\n    @(\"233\u00f76=38, 5\", \"375\u00f77=53, 4\"),\n    @(\"995\u00f74=248, 3\", \"775\u00f74=193, 3\"),\n    @(\"408\u00f78=51, 0\", \"227\u00f75=45, 2\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $oldText,    # FindText\n        $false,      # MatchCase\n        $true,       # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceOne)\n    )\n\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n"}
